# Auto commit at 2026-01-04  8:01:11.24
#
# Updates the Metrics sheet with new month/year/total figures, refreshes the
# "today" snapshot sheet (which pulls several of those figures in via
# formulas, plus a handful of its own static numbers), and leaves the
# workbook positioned/selected the way the author left it (today!C6 active,
# Metrics!D7 remembered as the old selection).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metrics sheet: fill in the previously-blank month cells and bump the
# running year/total figures to their new values.
# ---------------------------------------------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 51616.65
$metrics.Range("B3").Value = 39611.21
$metrics.Range("B4").Value = 13195.09
$metrics.Range("B5").Value = 2098

$metrics.Range("B6").Value = 5687487.3800000008
$metrics.Range("B7").Value = 4810328.84
$metrics.Range("B8").Value = 1677286.9100000001
$metrics.Range("B9").Value = 222375
$metrics.Range("B10").Value = 34152868.369999997
$metrics.Range("B11").Value = 32085604
$metrics.Range("B12").Value = 11959008.949999999
$metrics.Range("B13").Value = 1320005

# ---------------------------------------------------------------------
# "today" sheet: rows 15-22 get frozen ("paste values") to new numbers,
# losing their old =Metrics!Bxx formula, while rows 11-14 stay live
# formulas and simply recalc against the new Metrics values above.
# Rows 23-25 are independent static figures that also get refreshed.
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("today")

$today.Range("B15").Value = 5635870.7300000004
$today.Range("B16").Value = 4770717.63
$today.Range("B17").Value = 1664091.82
$today.Range("B18").Value = 220277
$today.Range("B19").Value = 34101251.719999999
$today.Range("B20").Value = 32045992.789999999
$today.Range("B21").Value = 11945813.859999999
$today.Range("B22").Value = 1317907

$today.Range("B23").Value = 50602.729999999996
$today.Range("B24").Value = 605487.01
$today.Range("B25").Value = 3256400.12

# ---------------------------------------------------------------------
# Restore the view/selection state: "today" becomes the active sheet with
# C6 selected, while Metrics keeps D7 remembered (but is no longer the
# active tab).
# ---------------------------------------------------------------------
$metrics.Activate()
$metrics.Range("D7").Select()

$today.Activate()
$today.Range("C6").Select()
